# Case and Fatality Demographics Data Updated
# Applies the 8/27/2021 refresh of the three "Fatalities by ..." sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Fatalities by Age Group
# ---------------------------------------------------------------------------
$wsAge = $wb.Worksheets.Item("Fatalities by Age Group")

$wsAge.Range("B2").Value  = 9
$wsAge.Range("B4").Value  = 41
$wsAge.Range("B5").Value  = 329
$wsAge.Range("B6").Value  = 1087
$wsAge.Range("B7").Value  = 3035
$wsAge.Range("B8").Value  = 6483
$wsAge.Range("B9").Value  = 5255
$wsAge.Range("B10").Value = 6634
$wsAge.Range("B11").Value = 7290
$wsAge.Range("B12").Value = 7141
$wsAge.Range("B13").Value = 17673

# ---------------------------------------------------------------------------
# Fatalities by Gender
# ---------------------------------------------------------------------------
$wsGender = $wb.Worksheets.Item("Fatalities by Gender")

$wsGender.Range("B2").Value = 23014
$wsGender.Range("B3").Value = 31976

# ---------------------------------------------------------------------------
# Fatalities by Race-Ethnicity
# ---------------------------------------------------------------------------
$wsRace = $wb.Worksheets.Item("Fatalities by Race-Ethnicity")

$wsRace.Range("B2").Value = 1130
$wsRace.Range("B3").Value = 5598
$wsRace.Range("B4").Value = 25108
$wsRace.Range("B5").Value = 302
$wsRace.Range("B6").Value = 22827
$wsRace.Range("B7").Value = 26

# ---------------------------------------------------------------------------
# Selections / active sheet. The last sheet activated/selected becomes the
# workbook's active tab, so activate "Fatalities by Race-Ethnicity" and
# "Fatalities by Gender" first, finishing on "Fatalities by Age Group" so it
# ends up as the active tab (matching the target workbookView.activeTab).
# ---------------------------------------------------------------------------
$wsRace.Activate()
$wsRace.Range("B8").Select()

$wsGender.Activate()
$wsGender.Range("B2:B4").Select()

$wsAge.Activate()
$wsAge.Range("F18").Select()
